# Refresh "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" worksheets, matching the newly generated
# gh-pages output snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new value, for the "展览" worksheet.
$exhibitionUpdates = @{
    "F2"  = 202
    "F3"  = 5506
    "F6"  = 29
    "F12" = 5017
    "F13" = 452
    "F15" = 198
    "F16" = 7
    "F18" = 4287
    "F22" = 53
    "F24" = 52
    "F25" = 157
    "F33" = 41
    "F34" = 42
}

# Row -> new value, for the "全部类型" worksheet (rows differ because this
# sheet aggregates events from every category sheet).
$allTypesUpdates = @{
    "F2"  = 202
    "F4"  = 5506
    "F7"  = 29
    "F13" = 5017
    "F14" = 452
    "F16" = 198
    "F17" = 7
    "F19" = 4287
    "F23" = 53
    "F25" = 52
    "F26" = 157
    "F34" = 41
    "F35" = 42
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($cellRef in $exhibitionUpdates.Keys) {
    $wsExhibition.Range($cellRef).Value = $exhibitionUpdates[$cellRef]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($cellRef in $allTypesUpdates.Keys) {
    $wsAllTypes.Range($cellRef).Value = $allTypesUpdates[$cellRef]
}
